$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed rows (player name, position, team) to reflect the new roster.
$ws.Range("A7").Value = "Jaxson Hayes"
$ws.Range("B7").Value = "PF,C"
$ws.Range("C7").Value = "Los Angeles Lakers"

$ws.Range("A9").Value = "Kevin Durant"
$ws.Range("B9").Value = "SF,PF"
$ws.Range("C9").Value = "Phoenix Suns"

$ws.Range("A13").Value = "Mark Williams"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Charlotte Hornets"

$ws.Range("A14").Value = "Karl-Anthony Towns"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "New York Knicks"

$ws.Range("A15").Value = "Tyrese Haliburton"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Indiana Pacers"

$ws.Range("A16").Value = "Cameron Johnson"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "Brooklyn Nets"

$ws.Range("A18").Value = "Klay Thompson"
$ws.Range("B18").Value = "SG,SF"
$ws.Range("C18").Value = "Dallas Mavericks"

# Remove the last row (Daniel Gafford / PF,C / Dallas Mavericks) entirely,
# shrinking the table from 18 to 17 data rows.
$ws.Rows.Item(19).Delete()

$wb.Save()
